$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.225455164909363
$ws.Range("B1").Value = 2.196571588516235
$ws.Range("C1").Value = 6.065976619720459
$ws.Range("D1").Value = 1.988811492919922
$ws.Range("E1").Value = 1.15572988986969
